$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-apply the same per-column font used by the existing rows (1-10) so the
# new A-column cells reuse the existing style slot (s="1") instead of
# minting a new one. (Column B's number format is applied per-cell, right
# after each value is set, below - doing it in bulk up front on still-blank
# cells would mint a throwaway style.)
$ws.Range("A11:A17").Font.Size = 10

# Row 11 - Joinville
$ws.Range("A11").Value = "Joinville"
$ws.Range("B11").Value = "'-26,3045"
$ws.Range("B11").NumberFormat = "#,##0"
$ws.Range("C11").Value = "'-48,8458"

# Row 12 - Lages
$ws.Range("A12").Value = "Lages"
$ws.Range("B12").Value = "'-27,815"
$ws.Range("B12").NumberFormat = "#,##0"
$ws.Range("C12").Value = "'-50,3259"

# Row 13 / Row 14 - Ouro Preto / Blumenal (entered out of simple row order,
# matching how the source workbook's shared strings ended up ordered)
$ws.Range("A13").Value = "Ouro Preto"
$ws.Range("C13").Value = "'-43,5031"
$ws.Range("A14").Value = "Blumenal"
$ws.Range("B14").Value = "'-26,9195"
$ws.Range("B14").NumberFormat = "#,##0"
$ws.Range("B13").Value = "'-20,3856"
$ws.Range("B13").NumberFormat = "#,##0"
$ws.Range("C14").Value = "'-49,066"

# Row 15 - Florianopolis
$ws.Range("A15").Value = "Florianopolis"
$ws.Range("B15").Value = "'-27,5949"
$ws.Range("B15").NumberFormat = "#,##0"
$ws.Range("C15").Value = "'-48,5482"

# Row 16 - Acre
$ws.Range("A16").Value = "Acre"
$ws.Range("B16").Value = "'-9,0238"
$ws.Range("B16").NumberFormat = "#,##0"
$ws.Range("C16").Value = "'-70,812"

# Row 17 - Fernando de Noronha
$ws.Range("A17").Value = "Fernando de Noronha"
$ws.Range("B17").Value = "'-3,8402"
$ws.Range("B17").NumberFormat = "#,##0"
$ws.Range("C17").Value = "'-32,4141"

# Leftover formatted-but-empty cells (mirrors the pre-existing E10 "ghost"
# formatting in the original sheet).
$ws.Range("D5").Font.Underline = 2
$ws.Range("A18").Font.Underline = 2

# Remove the underline formatting that used to be on C2.
$ws.Range("C2").Font.Underline = -4142
$ws.Range("C2").NumberFormat = "#,##0"
$ws.Range("C2").Value = "'-43,17"

# Window / view state tweaks captured in the source diff.
$excel.ActiveWindow.Zoom = 140
$ws.Range("C11").Select()
